$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row: "_old" suffix -> "_FV2210", "_new" suffix -> "_FV2304" ---
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    # columns A..J (1..10) -> "<name>_FV2210"
    $ws.Cells.Item(1, $i + 1).Value2 = $baseNames[$i] + "_FV2210"
    # columns L..U (12..21) -> "<name>_FV2304"
    $ws.Cells.Item(1, $i + 12).Value2 = $baseNames[$i] + "_FV2304"
}
# column K (11) stays "diff" - unchanged

# --- Turn the data range into an Excel Table (ListObject) with AutoFilter ---
$range = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"

# --- Freeze the header row (pane split after row 1) ---
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
